$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Text = "V roku Súhvezdie Herkules 2022: 13. – 22. júna, 12. – 21. júla, 10. – 19. augusta"
$find.Replacement.Text = "V roku 2022 môžete pozorovať súhvezdie Súhvezdie Herkules: 13. – 22. júna, 12. – 21. júla, 10. – 19. augusta"

$find.Execute(
    $find.Text,
    $false,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    $find.Replacement.Text,
    2
)
